$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# Bug fix: clear stale D:H numeric test values on row 7 and rows 29-47
# (these were left-over sample/debug numbers from a previous calibration
# run; graphs/labels pick up the blanks correctly once cleared).
$ws.Range("D7:H7").ClearContents()
$ws.Range("D29:H47").ClearContents()
